# Registration Data workbook update:
#  - split the "extra" registrants (rows 6-11) out of Sheet1 into a new Sheet2
#  - add Status/Remarks columns (H, I) to Sheet1 for the remaining registrants
#  - adjust views/selection/column widths to match

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Create Sheet2 (placed after Sheet1) and move the last six registrants
#    (old rows 6-11) there, keeping their formatting.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws1.Range("A6:G11").Copy()
$ws2.Range("A1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("A6:G11").Copy()
$ws2.Range("A1").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

# give the moved rows their (wrapped-text driven) row heights
$ws2.Rows.Item(1).RowHeight = 72.5
$ws2.Rows.Item(2).RowHeight = 72.5
$ws2.Rows.Item(3).RowHeight = 43.5
$ws2.Rows.Item(4).RowHeight = 72.5
$ws2.Rows.Item(5).RowHeight = 58
$ws2.Rows.Item(6).RowHeight = 58

# now remove the rows from Sheet1 (they now live on Sheet2)
$ws1.Range("A6:G11").Delete(-4162)   # xlShiftUp

# ---------------------------------------------------------------------------
# 2. Add the new "Status" / "Remarks" columns to the remaining Sheet1 rows.
# ---------------------------------------------------------------------------
$ws1.Range("H1").Value = "Status"
$ws1.Range("H1").Font.Bold = $true
$ws1.Range("H1").HorizontalAlignment = -4131

$ws1.Range("I1").Value = "Remarks"
$ws1.Range("I1").Font.Bold = $true
$ws1.Range("I1").HorizontalAlignment = -4131

$ws1.Range("H2").Value = "Pass"
$ws1.Range("I2").Value = "Registeration Success"

$ws1.Range("H3").Value = "Pass"
$ws1.Range("I3").Value = "Registeration Success"

$ws1.Range("H4").Value = "Pass"
$ws1.Range("I4").Value = "Registeration Success"

$ws1.Range("H5").Value = "Fail"
$ws1.Range("I5").Value = "Expected condition failed: waiting for presence of element located by: By.xpath: //p[text()='Let us know how we can help you! '] (tried for 10 second(s) with 500 milliseconds interval)"

$ws1.Columns.Item(9).ColumnWidth = 27.166666666666668

# ---------------------------------------------------------------------------
# 3. Views / selections.
# ---------------------------------------------------------------------------
$ws2.Range("N4").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C2").Select() | Out-Null
